$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "_old" -> "_FV2410", "_new" -> "_FV2504"
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Freeze the header row (pane split after row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

# 3. Turn the A1:U94 range into an Excel Table (adds autofilter + tableParts)
$range = $ws.Range("A1:U94")
$listObject = $ws.ListObjects.Add(1, $range, 0, 1)
$listObject.Name = "Table1"
